$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 33; this pushes the existing rows 33-66
# down to 34-67 and grows the used range from A1:R66 to A1:R67.
$ws.Rows.Item(33).Insert()

# Populate the newly inserted row 33 with the new weekly price record.
$ws.Range("A33").Value = 1
$ws.Range("B33").Value = "Agrícola del Norte S.A. de Arica"
$ws.Range("C33").Value = "Arica y Parinacota"
$ws.Range("D33").Value = 44669
$ws.Range("E33").Value = 15
$ws.Range("F33").Value = 100112021
$ws.Range("G33").Value = "Ají"
$ws.Range("H33").Value = "Inferno"
$ws.Range("I33").Value = "Primera"
$ws.Range("J33").Value = 120
$ws.Range("K33").Value = 28000
$ws.Range("L33").Value = 29000
$ws.Range("M33").Value = 28500
$ws.Range("N33").Value = "$/caja 15 kilos"
$ws.Range("O33").Value = "Región de Arica y Parinacota"
$ws.Range("P33").Value = 1900
$ws.Range("Q33").Value = 15
$ws.Range("R33").Value = "Hortaliza"
